{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph immediately preceding it) that followed the\n// last bibliography entry (\"VON SPERLING, M. ...\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\n// Find the \"Ver no Jupiter...\" paragraph; the empty paragraph right\n// before it and the \"\u00a9 2020...\" paragraph right after it are removed\n// together with it.\nlet verIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === targetTexts[0]) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex !== -1) {\n  const toDelete = [];\n  // Preceding blank paragraph (if present).\n  if (verIndex - 1 >= 0 && paragraphs.items[verIndex - 1].text.trim() === \"\") {\n    toDelete.push(paragraphs.items[verIndex - 1]);\n  }\n  // The \"Ver no Jupiter...\" paragraph itself.\n  toDelete.push(paragraphs.items[verIndex]);\n  // The following \"\u00a9 2020...\" paragraph (if present).\n  if (\n    verIndex + 1 < paragraphs.items.length &&\n    paragraphs.items[verIndex + 1].text.trim() === targetTexts[1]\n  ) {\n    toDelete.push(paragraphs.items[verIndex + 1]);\n  }\n\n  for (const para of toDelete) {\n    para.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph immediately preceding it) that followed the\n# last bibliography entry (\"VON SPERLING, M. ...\").\n\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"Contact: luizeleno@usp.br\"\n\n$verIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -ge 1) {\n    $indicesToDelete = New-Object System.Collections.ArrayList\n\n    # The following paragraph should mention the copyright/contact text.\n    if (($verIndex + 1) -le $d.Paragraphs.Count) {\n        $nextText = $d.Paragraphs.Item($verIndex + 1).Range.Text\n        if ($nextText -like \"*$copyrightText*\") {\n            [void]$indicesToDelete.Add($verIndex + 1)\n        }\n    }\n\n    [void]$indicesToDelete.Add($verIndex)\n\n    # The preceding blank paragraph (if present).\n    if ($verIndex - 1 -ge 1) {\n        $prevText = $d.Paragraphs.Item($verIndex - 1).Range.Text.Trim()\n        if ($prevText -eq \"\") {\n            [void]$indicesToDelete.Add($verIndex - 1)\n        }\n    }\n\n    # Delete from highest index to lowest so earlier indices stay valid.\n    $sorted = $indicesToDelete | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
